$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AY1").Value = "Transacción 06-2033"
$ws.Range("AZ1").Value = "Comentario de 06-2033"
$ws.Range("BA1").Value = "Transacción 06-2001"

$ws.Range("AY2").Value = "06-2033"
$ws.Range("AZ2").Value = "APROBADO"
$ws.Range("BA2").Value = "06-2001"
